$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Hydrogen): B3 value corrected, D3 cleared
$ws.Range("B3").Value = 798564.8337390764
$ws.Range("D3").Value = ""

# Row 4 (Methanol): C4 value corrected
$ws.Range("C4").Value = 28.84726982154575

# Row 5 (Ammonia): C5 value corrected
$ws.Range("C5").Value = 11692.64788055565

# Row 7: renamed from "Other" to "Biogas", D7 value corrected
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 1169.727739142278

# New row 8: "Other" (copy formatting from row 7's A column so the label
# style matches), with D8 holding the new value.
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = "Other"
$ws.Range("D8").Value = 753.6030022828552

$excel.CutCopyMode = $false
